$wb = $excel.ActiveWorkbook

# Belgium is worksheet #2 - use it as the template for the three new
# "market" sheets (Denmark, Sweden, Norway), just like UK/Belgium already
# are clones of the same layout.
$belgium = $wb.Worksheets.Item(2)

# --- Denmark ---------------------------------------------------------
$belgium.Copy($null, $belgium)
$denmark = $wb.Worksheets.Item(3)
$denmark.Name = "Denmark"
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Range("B4").Value = "NGC-3446/T2008"
$denmark.Activate()
$denmark.Range("A1:XFD1048576").Select()

# --- Sweden ------------------------------------------------------------
$denmark.Copy($null, $denmark)
$sweden = $wb.Worksheets.Item(4)
$sweden.Name = "Sweden"
$sweden.Range("B2").Value = "Sweden Market"
$sweden.Range("B4").Value = "NGC-3465/T2026"
$sweden.Activate()
$sweden.Range("A1:XFD1048576").Select()

# --- Norway --------------------------------------------------------
$sweden.Copy($null, $sweden)
$norway = $wb.Worksheets.Item(5)
$norway.Name = "Norway"
$norway.Range("B2").Value = "Norway Market"
$norway.Range("B4").Value = "NGC-3464/T1923"
$norway.Activate()
$norway.Range("B2:B4").Select()
$excel.ActiveCell = $norway.Range("B2")
